$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9783.27601853811
$ws.Range("F2").Value = -50.1346413531116

$ws.Range("C3").Value = 9780.7769109265
$ws.Range("F3").Value = 253.597292818899

$ws.Range("C4").Value = 9374.44961166726
$ws.Range("F4").Value = 269.83604185276

$ws.Range("C5").Value = 7205.84965290336
$ws.Range("F5").Value = 175.552776932418

$ws.Range("C6").Value = 7515.78854982245
$ws.Range("F6").Value = 199.82030333151

$ws.Range("C7").Value = 10858.6085828232
$ws.Range("F7").Value = 384.741371656185
